# Generate Report for Handoff
# Updates the localization status report: marks items as "Ready for handoff"
# (instead of "Handed back: in sync with en-US") and refreshes the
# handoff-generation timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Refreshed handoff-generation timestamps
$overview.Range("G2").Value = "2016-08-13 09:12:12"
$dede.Range("H2").Value = "2016-08-13 09:12:12"
$zhcn.Range("H2").Value = "2016-08-13 09:12:04"

# Narrow the Status/language columns now that the text is shorter
$overview.Range("E:E").ColumnWidth = 16.3
$overview.Range("F:F").ColumnWidth = 16.3
$zhcn.Range("C:C").ColumnWidth = 16.3
$dede.Range("C:C").ColumnWidth = 16.3
